$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.272.55"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -0.21%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.895.18"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +0.02%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.13%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "601.26"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.06%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "169.75"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +1.34%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.897.57"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.06%  "

$ws.Range("E8").Value = "  +0.03%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.530"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.22%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.165"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -1.79%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.42"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.58%  "

$ws.Range("E12").Value = "  -0.50%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000260"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +2.43%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.12"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.69%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.548.20"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.02%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.893.91"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.18%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "68.387.78"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.19%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.16"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +4.93%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.36"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -1.40%  "

$ws.Range("E20").Value = "  +0.17%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.78"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -2.69%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "469.02"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -4.48%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.739"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.66%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000164"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.93%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.50"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.50%  "

$ws.Range("E26").Value = "  +0.75%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.13"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +1.03%  "

$ws.Range("E28").Value = "  +0.02%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.98"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -1.75%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.96"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.82%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.045.45"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.00%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.87"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +1.86%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.31"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -2.40%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "31.35"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -1.28%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.43"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.50%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.865.60"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.42%  "

$ws.Range("E37").Value = "  -1.83%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.66"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +14.47%  "

$ws.Range("E39").Value = "  -0.29%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.140"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +0.66%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.92"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.47%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.999"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.11%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.313"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -1.18%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.000305"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +10.10%  "

$ws.Range("B45").Value = "Bittensor"
$ws.Range("C45").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "425.65"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.80%  "

$ws.Range("B46").Value = "Stacks"
$ws.Range("C46").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.98"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.12%  "

$ws.Range("E47").Value = "  -0.01%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.60"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.39%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "47.10"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -1.81%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "27.21"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +6.53%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "143.40"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.52%  "
